$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card4")

$ws.Range("A15").Value = "'4"
$ws.Range("L15").Value = "'11/11/2025"
$ws.Range("M15").Value = "قطع سير كويلر مسنن 1270"
$ws.Range("N15").Value = "تم تغير سير 1270(محمد نعيم)"
$ws.Range("O15").Value = "فني"
